$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Add the new hidden "DropdownOptions" sheet right after Sheet1, with the
#    percentage-bucket list used as the validation source.
# ---------------------------------------------------------------------------
$dropdownSheet = $wb.Worksheets.Add($null, $ws1, 1, $null)
$dropdownSheet.Name = "DropdownOptions"

$options = @("0% - 10%", "11% - 25%", "26% - 50%", "51% - 75%", "76% - 90%", "91% - 99%", "100%")
for ($i = 0; $i -lt $options.Length; $i++) {
    $dropdownSheet.Cells.Item($i + 1, 1).Value = $options[$i]
}

$dropdownSheet.Visible = $false

# ---------------------------------------------------------------------------
# 2. Fix the duplicated date/time number format (164 "yyyy-mm-dd h:mm:ss" vs
#    165 "YYYY-MM-DD HH:MM:SS") so the affected cells use the uppercase form.
# ---------------------------------------------------------------------------
$dateCells = @("T2", "V2", "W2", "P7", "Q7", "T7", "U7", "V7", "W7", "P8", "Q8", "T8", "U8", "V8", "W8")
foreach ($addr in $dateCells) {
    $ws1.Range($addr).NumberFormat = "yyyy-mm-dd h:mm:ss"
}

# ---------------------------------------------------------------------------
# 3. Clean up the stray empty "Actual Date of Completion" cells on rows 3-6.
# ---------------------------------------------------------------------------
foreach ($r in 3..6) {
    $ws1.Cells.Item($r, 17).ClearContents()
}

# ---------------------------------------------------------------------------
# 4. Add the new "Status as of July 4, 2025" column (AA) and its dropdown
#    data validation sourced from the hidden DropdownOptions sheet.
# ---------------------------------------------------------------------------
$ws1.Range("AA1").Value = "Status as of July 4, 2025"

$target = $ws1.Range("AA2:AA8")
$target.Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$target.Validation.IgnoreBlank = $true
$target.Validation.InCellDropdown = $true
$target.Validation.ShowInput = $false
$target.Validation.ShowError = $false

Write-Output "edit complete"
